# khurana_cv/data/cv_entries.xlsx — "working on extracting citations"
#
# Semantic content changes (everything else in the authoritative diff is
# just a shared-string re-index caused by deleting the now-unreferenced
# string "Awarded for highly-qualified incoming student from diverse
# background and experiences", which the engine garbage-collects on save
# once nothing points at it any more):
#
#   G15 (award / Lorraine Moe Davis Scholarship row): blank -> new note
#   G19 (award / Outstanding Poster, Grad Student Research Forum row): blank -> note
#   G20 (award / Winner, USAFacts Data Viz Competition row): note text tweaked
#   G21 (award / Graduate School Promising Scholar Award row): note text replaced
#
# G15 also picks up a brand-new 14pt "Source Sans Pro Reg" font run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 — additional_info was empty; give it new text in the larger font.
$ws.Range("G15").Value = "Demonstrated commitment to serving others"
$ws.Range("G15").Font.Size = 14
$ws.Range("G15").Font.Name = "Source Sans Pro Reg"

# Row 19 — additional_info was empty; now references the poster title.
$ws.Range("G19").Value = "For poster titled Diversity, Expenditure, and Achievement in US Public Schools"

# Row 20 — additional_info wording tweak ("poster" -> "final plots").
$ws.Range("G20").Value = "For final plots on Diversity, Expenditure, and Achievement in US Public Schools"

# Row 21 — additional_info replaced with the doctoral-recruitment description.
$ws.Range("G21").Value = "Recruit highly qualified incoming doctoral graduate students from diverse backgrounds who enhance and advance the academic and scholarly excellence of the university"
